$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 1081.7
$ws.Range("I4").Value = 1081.7
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1081.7
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -967.7
# Row 6
$ws.Range("H6").Value = 40226.434
$ws.Range("I6").Value = 112677.5
$ws.Range("K6").Value = 338032.5
$ws.Range("M6").Value = -337920.5
# Row 18
$ws.Range("H18").Value = 6118.125
$ws.Range("I18").Value = 7957.5
$ws.Range("K18").Value = 7957.5
$ws.Range("M18").Value = -7673.5
# Row 32
$ws.Range("H32").Value = 789.125
$ws.Range("J32").Value = 804.75
$ws.Range("L32").Value = 804.75
$ws.Range("N32").Value = -1456.75
# Row 33
$ws.Range("H33").Value = 613.65
$ws.Range("I33").Value = 259.25
$ws.Range("K33").Value = 259.25
$ws.Range("M33").Value = -30.25
# Row 62
$ws.Range("H62").Value = 142860590
$ws.Range("J62").Value = 3399
$ws.Range("L62").Value = 3399
$ws.Range("N62").Value = -4647
# Row 65
$ws.Range("H65").Value = 142860590
$ws.Range("J65").Value = 3399
$ws.Range("L65").Value = 16995
$ws.Range("N65").Value = -23235
# Row 100
$ws.Range("H100").Value = 12001552
$ws.Range("I100").Value = 19348046
$ws.Range("K100").Value = 19348046
$ws.Range("M100").Value = -19347505
# Row 107
$ws.Range("H107").Value = 4551.9565
$ws.Range("I107").Value = 4902.647
$ws.Range("J107").Value = 3558.3333
$ws.Range("K107").Value = 4902.647
$ws.Range("L107").Value = 3558.3333
$ws.Range("M107").Value = -2982.647
$ws.Range("N107").Value = -7398.3333
# Row 116
$ws.Range("H116").Value = 10437793
$ws.Range("I116").Value = 11386228
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 11386228
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -11382786
$ws.Range("N116").Value = -11884

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 15000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 15000
$ws.Range("K26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("M26").Value = 15000
$ws.Range("N26").Value = -15660
# Row 63
$ws.Range("H63").Value = 2169.9333
$ws.Range("I63").Value = 2169.9333
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2169.9333
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -1483.9333
# Row 66
$ws.Range("H66").Value = 2169.9333
$ws.Range("I66").Value = 2169.9333
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 10849.6665
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -7417.666500000001
# Row 74
$ws.Range("H74").Value = 7330.45
$ws.Range("I74").Value = 9369.923000000001
$ws.Range("J74").Value = 3542.8572
$ws.Range("K74").Value = 9369.923000000001
$ws.Range("L74").Value = 3542.8572
$ws.Range("M74").Value = -8495.923000000001
$ws.Range("N74").Value = -5290.8572
# Row 77
$ws.Range("H77").Value = 7330.45
$ws.Range("I77").Value = 9369.923000000001
$ws.Range("J77").Value = 3542.8572
$ws.Range("K77").Value = 46849.61500000001
$ws.Range("L77").Value = 17714.286
$ws.Range("M77").Value = -42481.61500000001
$ws.Range("N77").Value = -26450.286
# Row 110
$ws.Range("H110").Value = 2053.24
$ws.Range("I110").Value = 1666.7646
$ws.Range("J110").Value = 2874.5
$ws.Range("K110").Value = 1666.7646
$ws.Range("L110").Value = 2874.5
$ws.Range("M110").Value = 378.2354
$ws.Range("N110").Value = -6964.5
# Row 132
$ws.Range("H132").Value = 2448.0754
$ws.Range("I132").Value = 2107.8108
$ws.Range("K132").Value = 6323.432400000001
$ws.Range("M132").Value = -3793.432400000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3150.5715
$ws.Range("I20").Value = 1872.7142
$ws.Range("K20").Value = 1872.7142
$ws.Range("M20").Value = -1625.7142
# Row 105
$ws.Range("H105").Value = 131874
$ws.Range("I105").Value = 253498
$ws.Range("J105").Value = 10250
$ws.Range("K105").Value = 253498
$ws.Range("L105").Value = 10250
$ws.Range("M105").Value = -251751
$ws.Range("N105").Value = -13744
# Row 107
$ws.Range("H107").Value = 1951
$ws.Range("I107").Value = 2104.1667
$ws.Range("J107").Value = 113
$ws.Range("K107").Value = 2104.1667
$ws.Range("L107").Value = 113
$ws.Range("M107").Value = -184.1667000000002
$ws.Range("N107").Value = -3953

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9795.368
$ws.Range("I31").Value = 11153.071
$ws.Range("J31").Value = 5993.8
$ws.Range("K31").Value = 11153.071
$ws.Range("L31").Value = 5993.8
$ws.Range("M31").Value = -10858.071
$ws.Range("N31").Value = -6583.8
# Row 34
$ws.Range("H34").Value = 9795.368
$ws.Range("I34").Value = 11153.071
$ws.Range("J34").Value = 5993.8
$ws.Range("K34").Value = 11153.071
$ws.Range("L34").Value = 5993.8
$ws.Range("M34").Value = -10951.071
$ws.Range("N34").Value = -6397.8
# Row 134
$ws.Range("H134").Value = 8107.7
$ws.Range("I134").Value = 9943.733
$ws.Range("K134").Value = 29831.199
$ws.Range("M134").Value = -27296.199

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 850
$ws.Range("I3").Value = 850
$ws.Range("K3").Value = 2550
$ws.Range("M3").Value = -2438
# Row 4
$ws.Range("H4").Value = 70991620
$ws.Range("I4").Value = 71340056
$ws.Range("J4").Value = 70120530
$ws.Range("K4").Value = 214020168
$ws.Range("L4").Value = 210361590
$ws.Range("M4").Value = -214020056
$ws.Range("N4").Value = -210361814
# Row 11
$ws.Range("H11").Value = 29553.383
$ws.Range("I11").Value = 86.25
$ws.Range("K11").Value = 258.75
$ws.Range("M11").Value = -118.75
# Row 26
$ws.Range("H26").Value = 15.142858
$ws.Range("I26").Value = 4
$ws.Range("J26").Value = 56
$ws.Range("K26").Value = 12
$ws.Range("L26").Value = 168
$ws.Range("M26").Value = 276
$ws.Range("N26").Value = -744
# Row 29
$ws.Range("H29").Value = 320.6
$ws.Range("I29").Value = 227.95454
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 683.8636200000001
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = -406.8636200000001
$ws.Range("N29").Value = -3554
# Row 35
$ws.Range("H35").Value = 890.63635
$ws.Range("I35").Value = 793.375
$ws.Range("K35").Value = 2380.125
$ws.Range("M35").Value = -2092.125
# Row 40
$ws.Range("H40").Value = 59.1
$ws.Range("I40").Value = 32.857143
$ws.Range("J40").Value = 120.333336
$ws.Range("K40").Value = 131.428572
$ws.Range("L40").Value = 481.333344
$ws.Range("M40").Value = -62.428572
$ws.Range("N40").Value = -619.333344
# Row 47
$ws.Range("H47").Value = 646.5833
$ws.Range("I47").Value = 192
$ws.Range("J47").Value = 971.2857
$ws.Range("K47").Value = 576
$ws.Range("L47").Value = 2913.8571
$ws.Range("M47").Value = -145
$ws.Range("N47").Value = -3775.8571
# Row 55
$ws.Range("H55").Value = 11387.758
$ws.Range("I55").Value = 3002
$ws.Range("J55").Value = 11928.774
$ws.Range("K55").Value = 9006
$ws.Range("L55").Value = 35786.322
$ws.Range("M55").Value = -8829
$ws.Range("N55").Value = -36140.322
# Row 105
$ws.Range("H105").Value = 9950.625
$ws.Range("J105").Value = 9999.975
$ws.Range("L105").Value = 29999.925
$ws.Range("N105").Value = -35241.925
# Row 113
$ws.Range("H113").Value = 16599.857
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 16599.857
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").Value = 49799.571
$ws.Range("N113").Value = -54139.571
# Row 132
$ws.Range("H132").Value = 25131.143
$ws.Range("I132").Value = 1096.5714
$ws.Range("J132").Value = 37148.43
$ws.Range("K132").Value = 9869.142600000001
$ws.Range("L132").Value = 334335.87
$ws.Range("M132").Value = -7339.142600000001
$ws.Range("N132").Value = -339395.87

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7563.7915
$ws.Range("I70").Value = 6479.278
$ws.Range("J70").Value = 10817.333
$ws.Range("K70").Value = 6479.278
$ws.Range("L70").Value = 10817.333
$ws.Range("M70").Value = -6209.278
$ws.Range("N70").Value = -11357.333
# Row 73
$ws.Range("H73").Value = 7563.7915
$ws.Range("I73").Value = 6479.278
$ws.Range("J73").Value = 10817.333
$ws.Range("K73").Value = 6479.278
$ws.Range("L73").Value = 10817.333
$ws.Range("M73").Value = -5543.278
$ws.Range("N73").Value = -12689.333
# Row 102
$ws.Range("H102").Value = 6836.926
$ws.Range("I102").Value = 7666.95
$ws.Range("J102").Value = 4465.4287
$ws.Range("K102").Value = 7666.95
$ws.Range("L102").Value = 4465.4287
$ws.Range("M102").Value = -6044.95
$ws.Range("N102").Value = -7709.4287
# Row 126
$ws.Range("H126").Value = 11576.913
$ws.Range("I126").Value = 27342
$ws.Range("J126").Value = 7197.722
$ws.Range("K126").Value = 82026
$ws.Range("L126").Value = 21593.166
$ws.Range("M126").Value = -79556
$ws.Range("N126").Value = -26533.166

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
# Row 28
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
# Row 37
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
# Row 82
$ws.Range("H82").Value = 2763.0715
$ws.Range("I82").Value = 3023.889
$ws.Range("J82").Value = 2293.6
$ws.Range("K82").Value = 3023.889
$ws.Range("L82").Value = 2293.6
$ws.Range("M82").Value = -2662.889
$ws.Range("N82").Value = -3015.6
# Row 85
$ws.Range("H85").Value = 2763.0715
$ws.Range("I85").Value = 3023.889
$ws.Range("J85").Value = 2293.6
$ws.Range("K85").Value = 3023.889
$ws.Range("L85").Value = 2293.6
$ws.Range("M85").Value = -1775.889
$ws.Range("N85").Value = -4789.6
# Row 93
$ws.Range("H93").Value = 5064.05
$ws.Range("I93").Value = 7025.769
$ws.Range("K93").Value = 7025.769
$ws.Range("M93").Value = -5777.769

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Range("H74").Value = 154337.4
$ws.Range("J74").Value = 117796.75
$ws.Range("L74").Value = 117796.75
$ws.Range("N74").Value = -119668.75
# Row 77
$ws.Range("H77").Value = 154337.4
$ws.Range("J77").Value = 117796.75
$ws.Range("L77").Value = 353390.25
$ws.Range("N77").Value = -362750.25
